$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data corrections in the "uren registratie" (hours registration) block (rows 58-63)
$ws.Range("H58").Value = 0

$ws.Range("C59").Value = 4
$ws.Range("D59").Value = 4
$ws.Range("E59").Value = 4
$ws.Range("F59").Value = 4
$ws.Range("G59").Value = 4
$ws.Range("I59").Value = 4

# B63 is a hand-entered total (not a SUM formula like C63:I63)
$ws.Range("B63").Value = 8

# Restore the saved cursor/selection position
$ws.Range("K25").Select()
